$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数 / interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 43
$ws1.Range("F7").Value = 573
$ws1.Range("F8").Value = 19
$ws1.Range("F9").Value = 8236
$ws1.Range("F11").Value = 270
$ws1.Range("F12").Value = 1114
$ws1.Range("F13").Value = 833
$ws1.Range("F16").Value = 208
$ws1.Range("F17").Value = 103
$ws1.Range("F20").Value = 889

# Sheet "全部类型" (All types) - same underlying rows, shifted by one row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 43
$ws4.Range("F9").Value = 573
$ws4.Range("F10").Value = 19
$ws4.Range("F11").Value = 8236
$ws4.Range("F13").Value = 270
$ws4.Range("F14").Value = 1114
$ws4.Range("F15").Value = 833
$ws4.Range("F18").Value = 208
$ws4.Range("F19").Value = 103
$ws4.Range("F22").Value = 889
